# Mise à jour de l'application
# Adds a new training-session date column (AI) with attendance marks,
# mirroring the existing AH column's formatting/style, then updates
# the active selection like the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new date header (27/08/2025 -> Excel serial 45896) ---
$ws.Range("AI1").Value = 45896
$ws.Range("AH1").Copy()
$ws.Range("AI1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 2-27: attendance mark for the new session ---
# "P" = Présent, "R" = Réserve, "B" = Blessure, "RH" = Repos (row 27 is a
# constant "RH" banner row, not part of the player roster formulas).
$marks = @{
    2  = "P"
    3  = "R"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    12 = "P"
    13 = "P"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "B"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "P"
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "RH"
}

foreach ($r in 2..27) {
    $destCell = "AI" + $r
    $srcCell  = "AH" + $r

    # Set the value first, then copy the neighbouring cell's formatting so
    # the shared-formula ranges (which already span out to column VP/VQ/...)
    # recalculate immediately and pick up the correct number format/style.
    $ws.Range($destCell).Value = $marks[$r]
    $ws.Range($srcCell).Copy()
    $ws.Range($destCell).PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# Recalculate to make sure all dependent totals (COUNTA/COUNTIF columns) are fresh.
$excel.Calculate()

# --- Restore the author's last cursor position ---
[void]$ws.Range("AL25").Select()
